$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A7").Value = "Gabriel Pereira"
$ws.Range("B7").Value = 43902
$ws.Range("C7").Value = 0.70416666666666661
$ws.Range("D7").Value = "Bataille Navale"
$ws.Range("E7").Value = "Création du Menu Principale"
$ws.Range("F7").Value = "Créations du Menu Principale"

$ws.Range("A8").Value = "Gabriel Pereira"
$ws.Range("B8").Value = 43902
$ws.Range("C8").Value = 0.70416666666666661
$ws.Range("D8").Value = "Bataille Navale"
$ws.Range("E8").Value = "Création de la redirection"
$ws.Range("F8").Value = "Créations de la redirection"

$ws.Range("F8").Select()
